$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-28 all hold the date serial 45514 (2024-08-10)
# and need to be bumped by one day to 45515 (2024-08-11).
$ws.Range("C2:C28").Value = 45515
